$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.660.70"
$ws.Range("E2").Value = "  -5.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.201.62"
$ws.Range("E3").Value = "  -6.72%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.48"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.54"
$ws.Range("E6").Value = "  -8.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.587"
$ws.Range("E7").Value = "  -7.03%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.558"
$ws.Range("E9").Value = "  -8.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.65"
$ws.Range("E10").Value = "  -10.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.90"
$ws.Range("E11").Value = "  -3.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0825"
$ws.Range("E12").Value = "  -10.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.65"
$ws.Range("E13").Value = "  -9.71%  "
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.536.41"
$ws.Range("E15").Value = "  -6.60%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.853"
$ws.Range("E16").Value = "  -12.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.19"
$ws.Range("E17").Value = "  -6.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.206.93"
$ws.Range("E18").Value = "  -6.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.670.93"
$ws.Range("E19").Value = "  -5.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.44"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0958"
$ws.Range("E21").Value = "  -9.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.39"
$ws.Range("E22").Value = "  -11.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.02"
$ws.Range("E23").Value = "  -11.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.14"
$ws.Range("E24").Value = "  -10.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.45"
$ws.Range("E25").Value = "  -9.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.11"
$ws.Range("E26").Value = "  -8.60%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.05"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -9.94%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  -4.35%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.27"
$ws.Range("E31").Value = "  -12.68%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0882"
$ws.Range("E32").Value = "  -8.91%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.45"
$ws.Range("E33").Value = "  -8.40%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.99"
$ws.Range("E34").Value = "  -8.55%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.99"
$ws.Range("E35").Value = "  -7.61%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.80"
$ws.Range("E36").Value = "  -6.66%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.08"
$ws.Range("E37").Value = "  +7.71%  "
$ws.Range("E38").Value = "  -6.98%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("E39").Value = "  +7.17%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.108"
$ws.Range("E40").Value = "  -6.67%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.39"
$ws.Range("E41").Value = "  -5.94%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.79"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0323"
$ws.Range("E43").Value = "  -8.54%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.871.60"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.26"
$ws.Range("E46").Value = "  -4.55%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.34"
$ws.Range("E47").Value = "  -12.53%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.207"
$ws.Range("E48").Value = "  -9.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.32"
$ws.Range("E49").Value = "  -6.50%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.49"
$ws.Range("E50").Value = "  -13.15%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.01"
$ws.Range("E51").Value = "  -10.15%  "
